$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
$ws.Columns.Item(6).Delete()

$xlContinuous = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
$xlTop = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeTop
$xlBottom = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeBottom
$xlLeft = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeLeft
$xlRight = [Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight
$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
$xlLeftAlign = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignLeft

# Row 8 gets border formatting FIRST (closing row of the table), while style 7 (no-align,border8) is still used by B7/D7/E7
$r = $ws.Range("B8:E8")
$r.Borders.Item($xlTop).LineStyle = $xlContinuous
$r.Borders.Item($xlTop).Color = 0
$r.Borders.Item($xlTop).Weight = 2
$r.Borders.Item($xlBottom).LineStyle = $xlContinuous
$r.Borders.Item($xlBottom).Color = 0
$r.Borders.Item($xlBottom).Weight = 2

$a8 = $ws.Range("A8")
$a8.Borders.Item($xlTop).LineStyle = $xlContinuous
$a8.Borders.Item($xlTop).Color = 0
$a8.Borders.Item($xlTop).Weight = 2
$a8.Borders.Item($xlBottom).LineStyle = $xlContinuous
$a8.Borders.Item($xlBottom).Color = 0
$a8.Borders.Item($xlBottom).Weight = 2
$a8.Borders.Item($xlLeft).LineStyle = $xlContinuous
$a8.Borders.Item($xlLeft).Color = 0
$a8.Borders.Item($xlLeft).Weight = 2

$f8 = $ws.Range("F8")
$f8.Borders.Item($xlTop).LineStyle = $xlContinuous
$f8.Borders.Item($xlTop).Color = 0
$f8.Borders.Item($xlTop).Weight = 2
$f8.Borders.Item($xlBottom).LineStyle = $xlContinuous
$f8.Borders.Item($xlBottom).Color = 0
$f8.Borders.Item($xlBottom).Weight = 2
$f8.Borders.Item($xlRight).LineStyle = $xlContinuous
$f8.Borders.Item($xlRight).Color = 0
$f8.Borders.Item($xlRight).Weight = 2

# THEN row 7 alignment changes (frees up style 7 afterward)
$ws.Range("B7").HorizontalAlignment = $xlCenter
$ws.Range("C7:E7").HorizontalAlignment = $xlLeftAlign

Write-Host "done"
